$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.575.53"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").Value = "1.912.90"
$ws.Range("E3").Value = "  +5.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.52"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5091"
$ws.Range("E7").Value = "  +1.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3961"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09763"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("E10").Value = "  +5.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.24"
$ws.Range("E11").Value = "  +3.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.550"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.20"
$ws.Range("E13").Value = "  +3.95%  "
$ws.Range("D14").Value = "1.911.82"
$ws.Range("E14").Value = "  +5.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.591"
$ws.Range("E15").Value = "  +4.61%  "
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001140"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.02"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06667"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.19"
$ws.Range("E20").Value = "  +6.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.320"
$ws.Range("E22").Value = "  +7.15%  "
$ws.Range("D23").Value = "28.637.04"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.47"
$ws.Range("E24").Value = "  +3.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.283"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.735"
$ws.Range("E26").Value = "  +14.82%  "
$ws.Range("D27").Value = "2.138.52"
$ws.Range("E27").Value = "  +5.86%  "
$ws.Range("E28").Value = "  +4.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "159.72"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.73"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("E31").Value = "  +7.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1078"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.644"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.870"
$ws.Range("E35").Value = "  +11.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06817"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02446"
$ws.Range("E37").Value = "  +5.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.273"
$ws.Range("E38").Value = "  +9.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2234"
$ws.Range("E39").Value = "  +4.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.97"
$ws.Range("E40").Value = "  +6.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.107"
$ws.Range("E41").Value = "  +3.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6456"
$ws.Range("E42").Value = "  +4.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.190"
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.68"
$ws.Range("E45").Value = "  +3.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6102"
$ws.Range("E46").Value = "  +3.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.813"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.282"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.045"
$ws.Range("E49").Value = "  +5.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.34"
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("E51").Value = "  +3.24%  "
